$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.439.21"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.852.15"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'240.86"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07678"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").Value = "'0.2941"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'24.59"
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("D11").Value = "'0.07753"
$ws.Range("D12").Value = "1.854.49"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  +8.79%  "
$ws.Range("D14").Value = "'5.027"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "'0.6809"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "'83.61"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "2.101.82"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'6.164"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "29.463.03"
$ws.Range("D20").Value = "'229.27"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'7.453"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'156.87"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'0.1387"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").Value = "'8.400"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'1.319"
$ws.Range("E29").Value = "  +3.41%  "
$ws.Range("D30").Value = "'1.467"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "'0.05718"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").Value = "'4.134"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "'4.051"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'2.584"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "'2.782"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "'0.01794"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "1.218.20"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("D41").Value = "'6.509"
$ws.Range("E41").Value = "  +4.86%  "
$ws.Range("D42").Value = "'0.9075"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.010.73"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'101.80"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'66.55"
$ws.Range("E46").Value = "  +0.67%  "
$ws.Range("D47").Value = "'0.00000000120"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.135"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "'0.4019"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.985"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.684"
$ws.Range("E51").Value = "  -0.68%  "
